# Refactor element IDs - more accurate distribution component for avi 12,14,15,27,30
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: ehbv-avi.12-anseriformes-con -> ehbv-avi.12-anseridae-con
$ws.Range("A13").Value = "ehbv-avi.12-anseridae-con"
$ws.Range("B13").Value = "avi.12-anseridae"
$ws.Range("E13").Value = "ehbv-avi.12-anseridae-con"
$ws.Range("N13").Value = "anseridae"

# Row 15: ehbv-avi.14-gaviiformes -> ehbv-avi.14-gavia
$ws.Range("A15").Value = "ehbv-avi.14-gavia"
$ws.Range("B15").Value = "avi.14-gavia"
$ws.Range("E15").Value = "ehbv-avi.14-gavia"
$ws.Range("N15").Value = "gavia"

# Row 16: ehbv-avi.15-gaviiformes -> ehbv-avi.15-gavia
$ws.Range("A16").Value = "ehbv-avi.15-gavia"
$ws.Range("B16").Value = "avi.15-gavia"
$ws.Range("E16").Value = "ehbv-avi.15-gavia"
$ws.Range("N16").Value = "gavia"

# Row 27: ehbv-avi.27-suliformes-con -> ehbv-avi.27-sulidae-con
$ws.Range("A27").Value = "ehbv-avi.27-sulidae-con"
$ws.Range("B27").Value = "avi.27-sulidae"
$ws.Range("E27").Value = "ehbv-avi.27-sulidae-con"
$ws.Range("N27").Value = "sulidae"

# Row 30: ehbv-avi.30-anseriformes-con -> ehbv-avi.30-anseridae-con
$ws.Range("A30").Value = "ehbv-avi.30-anseridae-con"
$ws.Range("B30").Value = "avi.30-anseridae"
$ws.Range("E30").Value = "ehbv-avi.30-anseridae-con"
